# Update KHL referees stats workbook:
#  - "Главные" sheet (index 2): row 8 stat columns C:K change
#  - "Линейные" sheet (index 3): row 2 and row 14 stat columns C:K change
#  - all data rows (2-26) on both sheets get a refreshed as_of_utc timestamp (col AA)

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-03 14:31:22"

# --- Sheet "Главные" (2nd sheet) ---
$wsMain = $wb.Worksheets.Item(2)

# Row 8 updated stats
$wsMain.Range("C8").Value = 19
$wsMain.Range("D8").Value = 368
$wsMain.Range("E8").Value = 179
$wsMain.Range("F8").Value = 189
$wsMain.Range("G8").Value = 19.37
$wsMain.Range("H8").Value = 9.42
$wsMain.Range("I8").Value = 9.949999999999999
$wsMain.Range("J8").Value = 82
$wsMain.Range("K8").Value = 87

# Refresh as_of_utc timestamp for every data row
for ($r = 2; $r -le 26; $r++) {
    $wsMain.Range("AA" + $r).Value = $newTimestamp
}

# --- Sheet "Линейные" (3rd sheet) ---
$wsLine = $wb.Worksheets.Item(3)

# Row 2 updated stats
$wsLine.Range("C2").Value = 13
$wsLine.Range("D2").Value = 244
$wsLine.Range("E2").Value = 96
$wsLine.Range("F2").Value = 148
$wsLine.Range("G2").Value = 18.77
$wsLine.Range("H2").Value = 7.38
$wsLine.Range("I2").Value = 11.38
$wsLine.Range("J2").Value = 48
$wsLine.Range("K2").Value = 54

# Row 14 updated stats
$wsLine.Range("C14").Value = 21
$wsLine.Range("D14").Value = 348
$wsLine.Range("E14").Value = 174
$wsLine.Range("F14").Value = 174
$wsLine.Range("G14").Value = 16.57
$wsLine.Range("H14").Value = 8.289999999999999
$wsLine.Range("I14").Value = 8.289999999999999
$wsLine.Range("J14").Value = 87
$wsLine.Range("K14").Value = 82

# Refresh as_of_utc timestamp for every data row
for ($r = 2; $r -le 26; $r++) {
    $wsLine.Range("AA" + $r).Value = $newTimestamp
}
